$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 0.59215529571917125
$ws.Range("C2").Value = 1.6816769663972508
$ws.Range("D2").Value = 1.0600421034538978
$ws.Range("E2").Value = 1.5592843075931735

# Row 3 values
$ws.Range("B3").Value = 1.8279932993459795
$ws.Range("C3").Value = 1.4181226773019651
$ws.Range("D3").Value = 2.3739223718691895
$ws.Range("E3").Value = 1.2314511606815479

# Update the selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
